$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing text-typed "1" cell (C3) as a style/value template so the
# new "prio" cells stay shared-string text ("1") instead of becoming numeric 1.

# ---- Row 43: nudgincsm ----
$ws.Range("A43").Value = 'Eday'
$ws.Range("B43").Value = 'nudgincsm'
$ws.Range("C3").Copy($ws.Range("C43"))
$ws.Range("D43").Value = 'longitude latitude time'
$ws.Range("E43").Value = 'Nudging Increment of Water in Soil Moisture'
$ws.Range("F43").Value = 'kg m-2'
$ws.Range("G43").Formula = '=HYPERLINK("http://clipc-services.ceda.ac.uk/dreq/u/01c8c41a-a0d8-11e6-bc63-ac72891c3257.html","web")'
$ws.Range("H43").Value = 'To be implemented:  grib 126.151:  ifs code name = 151.126  part of MFPPHY.  Have to be  made available via PEXTRA, upto now with some  non-defined or adhoc grib code. Nudincsm is, consistent with sm, saved for each of the four soil layers'
$h43_1 = $ws.Range("H43").Characters(1, 19)
$h43_1.Font.Color = 0
$h43_1.Font.Name = "Calibri"
$h43_1.Font.Size = 11
$h43_20 = $ws.Range("H43").Characters(20, 13)
$h43_20.Font.Color = 1972430
$h43_20.Font.Name = "Calibri"
$h43_20.Font.Size = 11
$h43_33 = $ws.Range("H43").Characters(33, 207)
$h43_33.Font.Color = 0
$h43_33.Font.Name = "Calibri"
$h43_33.Font.Size = 11
$ws.Range("I43").Value = 'Emanuel Dutra, Wilhelm May, Thomas Reerink'
$ws.Range("J43").Value = 'A nudging increment refers to an amount added to parts of a model system. The phrase ''nudging_increment_in_X'' refers to an increment in quantity X over a time period which should be defined in the bounds of the time coordinate. ''Content'' indicates a quantity per unit area. ''Water'' means water in all phases. The mass content of water in soil refers to the vertical integral from the surface down to the bottom of the soil model. The ''soil content'' of a quantity refers to the vertical integral from the surface down to the bottom of the soil model. For the content between specified levels in the soil, standard names including ''content_of_soil_layer'' are used.'
$ws.Range("K43").Value = 'LS3MIP'

# ---- Row 44: nudgincswe ----
$ws.Range("A44").Value = 'Eday'
$ws.Range("B44").Value = 'nudgincswe'
$ws.Range("C3").Copy($ws.Range("C44"))
$ws.Range("D44").Value = 'longitude latitude time'
$ws.Range("E44").Value = 'Nudging Increment of Water in Snow'
$ws.Range("F44").Value = 'kg m-2'
$ws.Range("G44").Formula = '=HYPERLINK("http://clipc-services.ceda.ac.uk/dreq/u/0abbdddc-a0d8-11e6-bc63-ac72891c3257.html","web")'
$ws.Range("H44").Value = 'To be implemented:  grib 126.152:  ifs code name = 152.126  part of MFPPHY.  Have to be  made available via PEXTRA, upto now with some  non-defined or adhoc grib code.'
$h44_1 = $ws.Range("H44").Characters(1, 20)
$h44_1.Font.Color = 0
$h44_1.Font.Name = "Calibri"
$h44_1.Font.Size = 11
$h44_21 = $ws.Range("H44").Characters(21, 12)
$h44_21.Font.Color = 1972430
$h44_21.Font.Name = "Calibri"
$h44_21.Font.Size = 11
$h44_33 = $ws.Range("H44").Characters(33, 135)
$h44_33.Font.Color = 0
$h44_33.Font.Name = "Calibri"
$h44_33.Font.Size = 11
$ws.Range("I44").Value = 'Emanuel Dutra, Wilhelm May, Thomas Reerink'
$ws.Range("J44").Value = 'A nudging increment refers to an amount added to parts of a model system. The phrase ''nudging_increment_in_X'' refers to an increment in quantity X over a time period which should be defined in the bounds of the time coordinate. The surface called ''surface'' means the lower boundary of the atmosphere. ''Amount'' means mass per unit area. ''Snow and ice on land'' means ice in glaciers, ice caps, ice sheets & shelves, river and lake ice, any other ice on a land surface, such as frozen flood water, and snow lying on such ice or on the land surface.'
$ws.Range("K44").Value = 'LS3MIP'

# Match the recorded selection / scroll position from the edit session.
$ws.Range("B47").Select() | Out-Null
